# Braga.xlsx / "Rules" sheet: the single value that used to live in A1
# ("[Rule]") is moved down to E6, leaving A1 empty and extending the used
# range from A1:C3 to A1:E6. The selection ends up on B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Remember whichever sheet was active so we can restore it - changing the
# selection on $ws below must not change the workbook's active tab.
$originalActive = $wb.ActiveSheet

# Move the "[Rule]" label from A1 down to E6.
$ws.Range("A1").ClearContents()
$ws.Range("E6").Value = "[Rule]"

# Leave the selection on B8 of the Rules sheet, like in the source edit,
# then restore whatever sheet was originally active/selected.
$ws.Activate()
$ws.Range("B8").Select()
$originalActive.Activate()
